# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.986.35"
$ws.Range("E2").Value = "  +3.29%  "

# Row 3
$ws.Range("D3").Value = "3.056.29"
$ws.Range("E3").Value = "  +2.35%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.88"
$ws.Range("E5").Value = "  +2.75%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.58"
$ws.Range("E6").Value = "  +5.63%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").Value = "3.054.57"
$ws.Range("E8").Value = "  +2.23%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.508"
$ws.Range("E9").Value = "  +4.41%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("E10").Value = "  +6.27%  "

# Row 11
$ws.Range("E11").Value = "  -9.72%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.477"
$ws.Range("E12").Value = "  +7.46%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").Value = "  +6.04%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.00"
$ws.Range("E14").Value = "  +4.27%  "

# Row 15
$ws.Range("D15").Value = "3.542.60"
$ws.Range("E15").Value = "  +3.35%  "

# Row 16
$ws.Range("D16").Value = "64.015.13"
$ws.Range("E16").Value = "  +3.58%  "

# Row 17
$ws.Range("D17").Value = "3.054.97"
$ws.Range("E17").Value = "  +2.65%  "

# Row 18
$ws.Range("E18").Value = "  +1.82%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.72"
$ws.Range("E19").Value = "  +2.73%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.26"
$ws.Range("E20").Value = "  +2.75%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.03"
$ws.Range("E21").Value = "  +4.70%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.677"
$ws.Range("E22").Value = "  +4.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.57"
$ws.Range("E23").Value = "  +6.54%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.24"
$ws.Range("E24").Value = "  +13.58%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.56"
$ws.Range("E25").Value = "  +2.98%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.10%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.79"
$ws.Range("E27").Value = "  +3.18%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.92"
$ws.Range("E28").Value = "  +4.95%  "

# Row 29
$ws.Range("E29").Value = "  +2.28%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.26%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.23"
$ws.Range("E31").Value = "  +4.30%  "

# Row 32
$ws.Range("E32").Value = "  +1.66%  "

# Row 33
$ws.Range("E33").Value = "  +5.79%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.59"
$ws.Range("E34").Value = "  +2.09%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.19"
$ws.Range("E35").Value = "  +6.90%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.71"
$ws.Range("E36").Value = "  +1.19%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0405"
$ws.Range("E37").Value = "  +5.00%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "443.35"
$ws.Range("E38").Value = "  -0.79%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0806"
$ws.Range("E39").Value = "  +0.58%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.83"
$ws.Range("E40").Value = "  +15.56%  "

# Row 41
$ws.Range("D41").Value = "2.964.68"
$ws.Range("E41").Value = "  +1.20%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.21"
$ws.Range("E42").Value = "  +3.08%  "

# Row 43
$ws.Range("E43").Value = "  -0.87%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.62"
$ws.Range("E44").Value = "  +3.88%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.260"
$ws.Range("E45").Value = "  +5.68%  "

# Row 46
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.15"
$ws.Range("E46").Value = "  +8.12%  "

# Row 47
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  -0.02%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.113"
$ws.Range("E48").Value = "  +4.33%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.18"
$ws.Range("E49").Value = "  +2.99%  "

# Row 50
$ws.Range("D50").Value = "0.0₃0513"
$ws.Range("E50").Value = "  +5.92%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.07"
$ws.Range("E51").Value = "  +4.20%  "
